$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = [datetime]"2023-11-03"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
